# Uploads a new version of the day-sale / shortage report with an updated
# timestamp: one existing line (صوفي طويل جدا) picks up an extra sale
# (balance/turnover updated, sell price doubled because quantity sold
# doubled), and two new product lines are appended just above the
# totals/footer block. The totals row is shifted down accordingly and its
# total sell-price is bumped by the added amount, and the footer timestamp
# is refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $cellRef, $text) {
    # Force the value to be stored as text (shared string) even when it
    # looks like a pure number ("100.0000", "6.00", ...). We flip the
    # number format to Text, write the value, then restore the original
    # display format by pasting formats back from the cell itself so the
    # stored style id used on disk is unaffected.
    $style = $ws.Range($cellRef).NumberFormat
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------
# 1) Update the existing "صوفي طويل جدا" row (row 44): balance/turnover
#    and sell price reflect an additional sale during the day.
# ---------------------------------------------------------------------
$ws.Range("H44").Value = "0:0"
Set-TextCell $ws "P44" "100.0000"
$ws.Range("Q44").Value = "2:0"

# ---------------------------------------------------------------------
# 2) Make room for two new product rows right above the totals/footer
#    block: shift the totals row (old row 46) down to row 48 and the
#    footer row (old row 47) down to row 49, carrying their formatting
#    and values (and shared-string references) with them.
# ---------------------------------------------------------------------
$ws.Range("P46:Q46").UnMerge()
$ws.Range("A47:F47").UnMerge()
$ws.Range("G47:I47").UnMerge()
$ws.Range("K47:Q47").UnMerge()

# Footer row 47 -> row 49
$ws.Range("A47:Q47").Copy()
$ws.Range("A49:Q49").PasteSpecial(-4122) | Out-Null
$ws.Range("A47:Q47").Copy()
$ws.Range("A49:Q49").PasteSpecial(-4163) | Out-Null
$ws.Rows.Item(49).RowHeight = 16.5

# Totals row 46 -> row 48
$ws.Range("A46:Q46").Copy()
$ws.Range("A48:Q48").PasteSpecial(-4122) | Out-Null
$ws.Range("A46:Q46").Copy()
$ws.Range("A48:Q48").PasteSpecial(-4163) | Out-Null
$ws.Rows.Item(48).RowHeight = 24.75

# ---------------------------------------------------------------------
# 3) Populate the two freed-up rows (46, 47) with the new product lines,
#    cloning the normal data-row formatting from row 44.
# ---------------------------------------------------------------------
$ws.Range("A44:Q44").Copy()
$ws.Range("A46:Q46").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(46).RowHeight = 25.5

$ws.Range("A44:Q44").Copy()
$ws.Range("A47:Q47").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(47).RowHeight = 25.5

# Row 46: مناديل جيب مبلله
$ws.Range("A46").Value = 40
$ws.Range("C46").Value = "مناديل جيب مبلله "
$ws.Range("H46").Value = "11:0"
Set-TextCell $ws "L46" "0"
Set-TextCell $ws "N46" "6.00"
Set-TextCell $ws "P46" "0.0000"
$ws.Range("Q46").Value = "0:0"

# Row 47: مناديل مبلله كبيره
$ws.Range("A47").Value = 41
$ws.Range("C47").Value = "مناديل مبلله كبيره"
$ws.Range("H47").Value = "7:0"
Set-TextCell $ws "L47" "0"
Set-TextCell $ws "N47" "30.00"
Set-TextCell $ws "P47" "30.0000"
$ws.Range("Q47").Value = "1:0"

# ---------------------------------------------------------------------
# 4) Update the grand total and footer timestamp.
# ---------------------------------------------------------------------
$ws.Range("P48").Value = 1826.71
$ws.Range("A49").Value = "Monday, 28 July, 2025 4:22 PM"

# ---------------------------------------------------------------------
# 5) Re-establish merged cells for the shifted/new rows.
# ---------------------------------------------------------------------
$ws.Range("A46:B46").Merge()
$ws.Range("C46:G46").Merge()
$ws.Range("H46:K46").Merge()
$ws.Range("L46:M46").Merge()
$ws.Range("N46:O46").Merge()

$ws.Range("A47:B47").Merge()
$ws.Range("C47:G47").Merge()
$ws.Range("H47:K47").Merge()
$ws.Range("L47:M47").Merge()
$ws.Range("N47:O47").Merge()

$ws.Range("P48:Q48").Merge()
$ws.Range("A49:F49").Merge()
$ws.Range("G49:I49").Merge()
$ws.Range("K49:Q49").Merge()

Write-Output "edit complete"
